$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Update Coin name (B), Link (C), and Volume(1h) (E) for every row,
# plus Price (D) for rows whose new price contains two decimal separators
# (these can never be misread as a single number, so no special handling needed).

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '28.712.93'
$ws.Range('E2').Value = '  +2.69%  '
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '1.902.24'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('E4').Value = '  +3.24%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('E7').Value = '  +1.17%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('E8').Value = '  +3.68%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.896.88'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('B14').Value = 'BinanceUSD'
$ws.Range('C14').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E18').Value = '  +2.71%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '28.741.61'
$ws.Range('E22').Value = '  +2.65%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.131.42'
$ws.Range('E25').Value = '  +2.96%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E26').Value = '  +3.29%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E28').Value = '  -2.22%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E29').Value = '  +2.55%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E31').Value = '  +3.01%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E39').Value = '  +4.19%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E45').Value = '  +2.35%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E51').Value = '  +0.19%  '

# Step 2: Price (D) cells whose new value looks like a plain number would be
# auto-converted to a numeric cell (and lose e.g. trailing zeros) if we just set
# .Value directly, since the original sheet stores these as plain text. Mark the
# ranges as Text first, write the values, then drop the temporary number format so
# the cells end up as plain (unformatted) text cells again, matching the rest of the sheet.
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13:D21").NumberFormat = "@"
$ws.Range("D23:D24").NumberFormat = "@"
$ws.Range("D26:D51").NumberFormat = "@"
$ws.Range('D4').Value = '1.036'
$ws.Range('D5').Value = '319.85'
$ws.Range('D6').Value = '1.029'
$ws.Range('D7').Value = '0.5195'
$ws.Range('D8').Value = '0.3953'
$ws.Range('D9').Value = '0.08357'
$ws.Range('D10').Value = '1.134'
$ws.Range('D11').Value = '6.301'
$ws.Range('D13').Value = '20.55'
$ws.Range('D14').Value = '1.035'
$ws.Range('D15').Value = '7.308'
$ws.Range('D16').Value = '0.00001114'
$ws.Range('D17').Value = '91.69'
$ws.Range('D18').Value = '0.06822'
$ws.Range('D19').Value = '17.95'
$ws.Range('D20').Value = '1.032'
$ws.Range('D21').Value = '6.083'
$ws.Range('D23').Value = '11.24'
$ws.Range('D24').Value = '2.269'
$ws.Range('D26').Value = '162.75'
$ws.Range('D27').Value = '21.01'
$ws.Range('D28').Value = '2.450'
$ws.Range('D29').Value = '127.53'
$ws.Range('D30').Value = '0.1060'
$ws.Range('D31').Value = '1.056'
$ws.Range('D32').Value = '5.983'
$ws.Range('D33').Value = '3.681'
$ws.Range('D34').Value = '0.02468'
$ws.Range('D35').Value = '9.408'
$ws.Range('D36').Value = '0.06634'
$ws.Range('D37').Value = '0.2211'
$ws.Range('D38').Value = '0.6560'
$ws.Range('D39').Value = '1.258'
$ws.Range('D40').Value = '1.191'
$ws.Range('D41').Value = '5.016'
$ws.Range('D42').Value = '11.20'
$ws.Range('D43').Value = '0.6142'
$ws.Range('D44').Value = '13.33'
$ws.Range('D45').Value = '3.761'
$ws.Range('D46').Value = '1.303'
$ws.Range('D47').Value = '2.023'
$ws.Range('D48').Value = '1.232'
$ws.Range('D49').Value = '123.08'
$ws.Range('D50').Value = '0.06985'
$ws.Range('D51').Value = '78.06'
$ws.Range("D4:D11").ClearFormats()
$ws.Range("D13:D21").ClearFormats()
$ws.Range("D23:D24").ClearFormats()
$ws.Range("D26:D51").ClearFormats()
